$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "59.004.99"
$ws.Range("E2").Value = "  +0.76%  "
$ws.Range("D3").Value = "2.496.29"
$ws.Range("E3").Value = "  +1.42%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").Value = "'533.97"
$ws.Range("E5").Value = "  +0.50%  "
$ws.Range("D6").Value = "'136.05"
$ws.Range("E6").Value = "  +0.82%  "
$ws.Range("D7").Value = "'0.998"
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("D8").Value = "'0.566"
$ws.Range("E8").Value = "  +1.44%  "
$ws.Range("D9").Value = "'0.102"
$ws.Range("E9").Value = "  +2.49%  "
$ws.Range("E10").Value = "  -1.28%  "
$ws.Range("D11").Value = "'5.40"
$ws.Range("E11").Value = "  +1.85%  "
$ws.Range("D12").Value = "'0.348"
$ws.Range("E12").Value = "  +1.30%  "
$ws.Range("D13").Value = "2.943.60"
$ws.Range("E13").Value = "  +1.49%  "
$ws.Range("D14").Value = "58.905.19"
$ws.Range("E14").Value = "  +0.70%  "
$ws.Range("D15").Value = "'22.66"
$ws.Range("E15").Value = "  -0.78%  "
$ws.Range("E16").Value = "  +0.22%  "
$ws.Range("D17").Value = "2.502.48"
$ws.Range("E17").Value = "  +1.30%  "
$ws.Range("D18").Value = "'11.04"
$ws.Range("E18").Value = "  +1.92%  "
$ws.Range("D19").Value = "'4.25"
$ws.Range("E19").Value = "  +1.14%  "
$ws.Range("D20").Value = "'323.27"
$ws.Range("E20").Value = "  +0.14%  "
$ws.Range("B21").Value = "Uniswap"
$ws.Range("C21").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D21").Value = "'5.99"
$ws.Range("E21").Value = "  +4.14%  "
$ws.Range("B22").Value = "Dai"
$ws.Range("C22").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D22").Value = "'1.00"
$ws.Range("E22").Value = "  +0.26%  "
$ws.Range("D23").Value = "'65.21"
$ws.Range("E23").Value = "  +3.94%  "
$ws.Range("D24").Value = "'0.421"
$ws.Range("E24").Value = "  +2.14%  "
$ws.Range("E25").Value = "  -0.06%  "
$ws.Range("E26").Value = "  +1.34%  "
$ws.Range("D27").Value = "'7.53"
$ws.Range("E27").Value = "  +1.05%  "
$ws.Range("D28").Value = "0.0₃0761"
$ws.Range("E28").Value = "  +0.78%  "
$ws.Range("B29").Value = "Aptos"
$ws.Range("C29").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D29").Value = "'6.48"
$ws.Range("E29").Value = "  -3.08%  "
$ws.Range("B30").Value = "Monero"
$ws.Range("C30").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D30").Value = "'170.49"
$ws.Range("E30").Value = "  +3.08%  "
$ws.Range("B31").Value = "PancakeSwap"
$ws.Range("C31").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D31").Value = "'1.74"
$ws.Range("E31").Value = "  -0.57%  "
$ws.Range("D32").Value = "'1.18"
$ws.Range("E32").Value = "  +5.00%  "
$ws.Range("D33").Value = "'0.998"
$ws.Range("E33").Value = "  +0.00%  "
$ws.Range("D34").Value = "'18.36"
$ws.Range("E34").Value = "  +0.39%  "
$ws.Range("D35").Value = "'1.35"
$ws.Range("E35").Value = "  -0.64%  "
$ws.Range("D36").Value = "'4.04"
$ws.Range("E36").Value = "  +0.04%  "
$ws.Range("D37").Value = "'1.53"
$ws.Range("E37").Value = "  -0.58%  "
$ws.Range("D38").Value = "'3.57"
$ws.Range("E38").Value = "  +0.39%  "
$ws.Range("B39").Value = "SuiNetwork"
$ws.Range("C39").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D39").Value = "'0.799"
$ws.Range("E39").Value = "  -0.46%  "
$ws.Range("B40").Value = "Bittensor"
$ws.Range("C40").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D40").Value = "'282.26"
$ws.Range("E40").Value = "  +2.16%  "
$ws.Range("D41").Value = "'0.998"
$ws.Range("D42").Value = "'5.03"
$ws.Range("E42").Value = "  -1.65%  "
$ws.Range("B43").Value = "Aave"
$ws.Range("C43").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D43").Value = "'130.10"
$ws.Range("E43").Value = "  +5.82%  "
$ws.Range("B44").Value = "Mantle"
$ws.Range("C44").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D44").Value = "'0.602"
$ws.Range("E44").Value = "  +2.16%  "
$ws.Range("D45").Value = "'10.92"
$ws.Range("E45").Value = "  +0.63%  "
$ws.Range("D46").Value = "'0.0924"
$ws.Range("E46").Value = "  +0.08%  "
$ws.Range("D47").Value = "'0.0499"
$ws.Range("E47").Value = "  -1.01%  "
$ws.Range("D48").Value = "'0.0218"
$ws.Range("E48").Value = "  -0.19%  "
$ws.Range("D49").Value = "'17.27"
$ws.Range("E49").Value = "  +0.23%  "
$ws.Range("D50").Value = "1.756.29"
$ws.Range("E50").Value = "  +0.40%  "
$ws.Range("D51").Value = "'0.982"
$ws.Range("E51").Value = "  +0.25%  "
